# Germany Landesliga update (28-04-2024 15:37)
# Rows 4 and 5 (match ids 2 and 3, played on the same date) had their
# match data (everything except id/Div/Date) swapped between the two
# rows, so the "SV Schott Jena" vs "SSV Markranstadt" fixture entries
# end up pointing at the correct match results / odds.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

foreach ($col in $cols) {
    $addr4 = "$col" + "4"
    $addr5 = "$col" + "5"
    $val4 = $ws.Range($addr4).Value2
    $val5 = $ws.Range($addr5).Value2
    $ws.Range($addr4).Value = $val5
    $ws.Range($addr5).Value = $val4
}
